$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new record row (row 20) with the latest fly detection data.
$row = 20

$ws.Range("A$row").Value = "89ec9c17-dac0-435a-851d-754b073844e5"
$ws.Range("B$row").Value = "mosca"

$ws.Range("C$row").Value = 45891
$ws.Range("C$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("D$row").Value = "image_20250822211929_ppp0.jpg"
$ws.Range("E$row").Value = "PLACA_20250717165933"
$ws.Range("F$row").Value = "Beja"
$ws.Range("G$row").Value = 38.02035
$ws.Range("H$row").Value = -7.94715
$ws.Range("I$row").Value = "1022,1004,1060,1059"

$ws.Range("J$row").Value = "'0.68"
$ws.Range("J$row").Style = "Normal"
